# Updated cryptos list on Mon Jul 29 20:36:03 UTC 2024 with GitHub Actions
#
# Applies refreshed Price / Volume(1h) figures (and a handful of row
# re-rankings where a coin entered/left/moved within the top list) to the
# "cryptos" worksheet.
#
# NOTE: every value in columns B:E of this sheet is stored as plain text in
# the source data (coin names, coinranking.com URLs, price strings that can
# contain more than one '.' as a thousands separator, and percentage
# strings padded with spaces). Excel's COM layer normally "helpfully"
# re-interprets a numeric-looking string typed into `.Value` as a real
# number (stripping formatting like trailing zeros, e.g. "2.40" -> 2.4).
# To faithfully reproduce the original text values we force each target
# cell to Text format before writing, then restore the default "Normal"
# style so we don't leave a stray text-format style applied to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "67.237.80"
Set-TextValue "E2" "  -1.14%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.325.72"
Set-TextValue "E3" "  +1.85%  "

# Row 5 - Solana
Set-TextValue "D5" "186.56"
Set-TextValue "E5" "  +1.04%  "

# Row 6 - BNB
Set-TextValue "D6" "578.13"
Set-TextValue "E6" "  -0.88%  "

# Row 7 - XRP
Set-TextValue "E7" "  +1.16%  "

# Row 8 - USDC
Set-TextValue "E8" "  -0.03%  "

# Row 9 - Dogecoin
Set-TextValue "D9" "0.129"
Set-TextValue "E9" "  -0.06%  "

# Row 10 - Toncoin
Set-TextValue "E10" "  +0.71%  "

# Row 11 - Cardano
Set-TextValue "E11" "  +0.11%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "3.893.21"
Set-TextValue "E12" "  +1.47%  "

# Row 13 - TRON
Set-TextValue "E13" "  -0.59%  "

# Row 14 - Avalanche
Set-TextValue "D14" "27.44"
Set-TextValue "E14" "  -0.03%  "

# Row 15 - WrappedBTC
Set-TextValue "D15" "67.452.43"

# Row 16 - ShibaInu
Set-TextValue "E16" "  +0.13%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "3.307.85"
Set-TextValue "E17" "  +2.40%  "

# Row 18 - BitcoinCash
Set-TextValue "D18" "445.07"
Set-TextValue "E18" "  +6.94%  "

# Row 19 - Polkadot
Set-TextValue "E19" "  -0.47%  "

# Row 20 - Chainlink
Set-TextValue "D20" "13.59"
Set-TextValue "E20" "  +2.21%  "

# Row 21 - Uniswap
Set-TextValue "D21" "7.74"
Set-TextValue "E21" "  +2.78%  "

# Row 22 - Litecoin
Set-TextValue "D22" "74.12"
Set-TextValue "E22" "  +3.74%  "

# Row 23 - Dai
Set-TextValue "E23" "  -0.02%  "

# Row 24 - Polygon
Set-TextValue "E24" "  +1.87%  "

# Row 25 - WrappedeETH
Set-TextValue "D25" "3.461.72"
Set-TextValue "E25" "  +1.51%  "

# Row 26 - PEPE
Set-TextValue "E26" "  +1.79%  "

# Row 27 - Kaspa
Set-TextValue "E27" "  +1.19%  "

# Row 28 - InternetComputer(DFINITY)
Set-TextValue "D28" "9.07"
Set-TextValue "E28" "  -3.84%  "

# Row 29 - Binance-PegBSC-USD
Set-TextValue "D29" "0.994"
Set-TextValue "E29" "  -0.77%  "

# Row 30 - PancakeSwap
Set-TextValue "E30" "  +1.35%  "

# Row 31 - EthereumClassic
Set-TextValue "D31" "22.96"
Set-TextValue "E31" "  +1.35%  "

# Row 32 - NEARProtocol
Set-TextValue "E32" "  -2.03%  "

# Row 33 - USDe
Set-TextValue "D33" "0.998"
Set-TextValue "E33" "  -0.02%  "

# Row 34 - Aptos
Set-TextValue "E34" "  -0.28%  "

# Row 35 - Fetch.AI
Set-TextValue "E35" "  -0.20%  "

# Row 36 - ImmutableX
Set-TextValue "E36" "  +5.70%  "

# Row 37 - Monero
Set-TextValue "D37" "162.66"
Set-TextValue "E37" "  +0.02%  "

# Row 38 - was Stacks, now EnergySwap
Set-TextValue "B38" "EnergySwap"
Set-TextValue "C38" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D38" "27.42"
Set-TextValue "E38" "  +1.53%  "

# Row 39 - was EnergySwap, now Stacks
Set-TextValue "B39" "Stacks"
Set-TextValue "C39" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D39" "1.85"
Set-TextValue "E39" "  -1.64%  "

# Row 40 - Maker
Set-TextValue "D40" "2.785.28"
Set-TextValue "E40" "  +5.54%  "

# Row 41 - Mantle
Set-TextValue "D41" "0.792"
Set-TextValue "E41" "  -0.65%  "

# Row 42 - Filecoin
Set-TextValue "E42" "  +0.55%  "

# Row 43 - RenderToken
Set-TextValue "D43" "6.24"
Set-TextValue "E43" "  -1.77%  "

# Row 44 - was InjectiveProtocol, now Hedera
Set-TextValue "B44" "Hedera"
Set-TextValue "C44" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D44" "0.0672"
Set-TextValue "E44" "  -0.40%  "

# Row 45 - was Hedera, now InjectiveProtocol
Set-TextValue "B45" "InjectiveProtocol"
Set-TextValue "C45" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D45" "24.82"
Set-TextValue "E45" "  +2.07%  "

# Row 46 - OKB
Set-TextValue "D46" "40.15"
Set-TextValue "E46" "  -1.55%  "

# Row 47 - dogwifhat
Set-TextValue "D47" "2.40"
Set-TextValue "E47" "  -0.87%  "

# Row 48 - Bittensor
Set-TextValue "D48" "326.70"
Set-TextValue "E48" "  -3.07%  "

# Row 49 - VeChain
Set-TextValue "E49" "  +0.06%  "

# Row 50 - ONDO
Set-TextValue "D50" "0.990"
Set-TextValue "E50" "  +1.44%  "

# Row 51 - was Arweave, now Cosmos
Set-TextValue "B51" "Cosmos"
Set-TextValue "C51" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D51" "6.20"
Set-TextValue "E51" "  -1.07%  "
